$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-07-20"

# Update the row label for July to reflect the new "through" date
$ws.Range("A8").Value = "July (through 07-20)"

# Update July row (row 8) values for each year column (B..I)
$ws.Range("B8").Value = 26
$ws.Range("C8").Value = 41
$ws.Range("D8").Value = 42
$ws.Range("E8").Value = 51
$ws.Range("F8").Value = 31
$ws.Range("G8").Value = 82
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 112

# Update Total row (row 9) values for each year column (B..I)
$ws.Range("B9").Value = 151
$ws.Range("C9").Value = 289
$ws.Range("D9").Value = 432
$ws.Range("E9").Value = 404
$ws.Range("F9").Value = 282
$ws.Range("G9").Value = 554
$ws.Range("H9").Value = 860
$ws.Range("I9").Value = 918
